# A new weekly price record was inserted as row 51 (pushing the existing
# rows 51-78 down to 52-79). The sheet's used range grows from A1:R78 to
# A1:R79 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51; everything below shifts down one row.
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new weekly record.
$ws.Cells.Item(51, 1).Value = 6
$ws.Cells.Item(51, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(51, 3).Value = "Metropolitana"
$ws.Cells.Item(51, 4).Value = 45086
$ws.Cells.Item(51, 5).Value = 13
$ws.Cells.Item(51, 6).Value = 100112035
$ws.Cells.Item(51, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 400
$ws.Cells.Item(51, 11).Value = 16000
$ws.Cells.Item(51, 12).Value = 18000
$ws.Cells.Item(51, 13).Value = 16850
$ws.Cells.Item(51, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(51, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(51, 16).Value = 1123
$ws.Cells.Item(51, 17).Value = 15
$ws.Cells.Item(51, 18).Value = "Hortaliza"

# Keep the date column's formatting consistent with the rest of column D.
$ws.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
